$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "NA" row (row 4) from the education options list, shifting
# the rows below it up by one — mirrors selecting row 4 and deleting it.
$ws.Rows.Item(4).Delete()

# Leave the selection on the row that now occupies the deleted row's
# position, matching the post-delete selection state.
$ws.Range("A4:XFD4").Select()
